$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Pass"
$ws.Range("D3").Value = "Pass"
